# "Generate Report for Archive"
# Updates the localization status text from "Ready for handoff" to "In Translation"
# on all three sheets, and shrinks the now-narrower "Status" column widths
# (Overview!E:F, zh-cn!C, de-de!C) to match the regenerated report.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# The runtime quantizes ColumnWidth to the nearest 1/6 character-width unit
# when persisting to OOXML; 12.5 is the input that lands on the stored width
# closest to the target 13.4101848602295 (i.e. 13.333333333333334).
$newWidth  = 12.5

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: column C ("Status") holds the status text ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: column C ("Status") holds the status text ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
